# Add a "PMID" column to the "studies" sheet and a "notes" column to the
# "counts" sheet, then leave the "counts" sheet as the active/selected one
# (matching the author's recorded cursor positions).

$wb = $excel.ActiveWorkbook

# --- studies sheet: new column H = "PMID" ---------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Cells.Item(1, 8).Value = "PMID"
[void]$wsStudies.Range("H2").Select()

# --- counts sheet: new column F = "notes" ----------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Cells.Item(1, 6).Value = "notes"

# Make "counts" the active sheet/tab, with F2 selected.
[void]$wsCounts.Activate()
[void]$wsCounts.Range("F2").Select()
